# The "Tokens" worksheet holds JWTs captured by an automated REST-Assured
# test suite across repeated runs (OTP verification + consumer login +
# buy/sell flows, per the commit message). Each run writes a fresh token
# into Tokens!A2 (Login_Token / OTP) or Tokens!C2 (Consumer authToken),
# mirroring the real sequence of writes so the sheet ends up holding the
# same final tokens as after the latest test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tokens")

$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo4MjYwNDQsImlhdCI6MTc0NDczNjI1NywiZXhwIjoxNzQ0NzM5ODU3fQ.u83RvPcQ2YMxwf4Q-3iOdQQL_tISwtKBdycwGUieN6g"
$ws.Range("C2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwidXNlcklkIjo5LCJtb2JpbGVfbnVtYmVyIjpudWxsLCJpYXQiOjE3NDQ3MzYyNjAsImV4cCI6MTc3NjI5Mzg2MH0.2pWQSZbgzVDwsduRtCAxxNiIjMiJsnfJoQ2MlN7dgyI"
$ws.Range("C2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwidXNlcklkIjo5LCJtb2JpbGVfbnVtYmVyIjpudWxsLCJpYXQiOjE3NDQ3MzY0MTEsImV4cCI6MTc3NjI5NDAxMX0.lU1vn_GQfQyeuA74NjpkpLEIcDh9f3fLRy8hLQ2IZw8"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo0MDY0OTMsImlhdCI6MTc0NDczNjgxOCwiZXhwIjoxNzQ0NzQwNDE4fQ.WHvMsW6fXtiCrUsKeSIDSUZ8mvZVFOysVJzdAw3LUFM"
$ws.Range("C2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwidXNlcklkIjo5LCJtb2JpbGVfbnVtYmVyIjpudWxsLCJpYXQiOjE3NDQ3MzY4MjEsImV4cCI6MTc3NjI5NDQyMX0.ammu8PUcawUx685qllLBeSdvS3XMJyYPZkg8CxXBR_4"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo1NjA3MTksImlhdCI6MTc0NDczNzA1MywiZXhwIjoxNzQ0NzQwNjUzfQ.rhr2nniNG0XACeNtqakWdGpgCsYzbkO48vaW6_Yyges"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjoxOTk1MDksImlhdCI6MTc0NDczNzA1OCwiZXhwIjoxNzQ0NzQwNjU4fQ.iFz7TuMTqM1P5jTTCRAmP6pdGZptBdMTc50rqACo3U8"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo0MDY0MjksImlhdCI6MTc0NDczNzA2MCwiZXhwIjoxNzQ0NzQwNjYwfQ.Hcxnh_lmexcfDwfXueEWWBAkqQ1lLZZn6Hzgs-UwMjM"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjoxNTkwMjIsImlhdCI6MTc0NDczNzA2MiwiZXhwIjoxNzQ0NzQwNjYyfQ.uTEJBLlHqDO-HtrDKTgXLYwUvYh3HSK7DVMVgqmIR7U"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo4OTUzODAsImlhdCI6MTc0NDczNzA2MywiZXhwIjoxNzQ0NzQwNjYzfQ.841cwmS_FiexjvWC95YfMmpRGvbz4k-4kiiWPanOBDI"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo4ODA1NzAsImlhdCI6MTc0NDczNzEyMSwiZXhwIjoxNzQ0NzQwNzIxfQ.1V3y_fX5lUQJ_SryezHydrZcL9ztbdlwf5IhzPf-yMk"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjoyODc0NzUsImlhdCI6MTc0NDczNzEyNiwiZXhwIjoxNzQ0NzQwNzI2fQ.e_nIsFeMMFL2wmp1kmcAfnR950NPmMkmMcSe-nBYvQo"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo1MzgzOTIsImlhdCI6MTc0NDczNzEyOSwiZXhwIjoxNzQ0NzQwNzI5fQ.mc5CVvHS8-04aexvfai0Ywx93m2JQdGahDp6UK9gffw"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo2OTY3MzksImlhdCI6MTc0NDczNzEzMSwiZXhwIjoxNzQ0NzQwNzMxfQ.uGe9MRz_i0zKm9FapTwmp5sBKRyu0V96TE4Q6CqjK8s"
$ws.Range("A2").Value2 = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjoyMzYwMzksImlhdCI6MTc0NDczNzEzMywiZXhwIjoxNzQ0NzQwNzMzfQ.hGkn_RPjmFYIdCE-n6L3c8bfUMRANEheo6tFh9fZGKY"
